$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Collapse the three long "CORE COMPETENCIES" bullet paragraphs
#    into a single short paragraph with just the three category
#    names separated by a bullet character.
# -----------------------------------------------------------------
$bullet = [char]0x2022

$skillsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Statistical Analysis & Machine Learning: Advanced Statistical Modeling")) {
        $skillsPara = $p
        break
    }
}
if ($skillsPara -eq $null) {
    throw "Could not locate the CORE COMPETENCIES 'Statistical Analysis & Machine Learning' paragraph"
}

$skillsPara.Range.Text = "Statistical Analysis & Machine Learning " + $bullet + " Big Data & Data Engineering " + $bullet + " Data Visualization & Reporting"

# The next two paragraphs (Big Data & Data Engineering / Data
# Visualization & Reporting long-form text) are now obsolete - drop them.
$next1 = $skillsPara.Next()
$next1.Range.Delete()
$next2 = $skillsPara.Next()
$next2.Range.Delete()

# -----------------------------------------------------------------
# 2) Append a new "TECHNICAL SKILLS" section at the end of the
#    document, with the same three categories but written out as
#    their own paragraphs (Heading2 for the title, Normal body
#    paragraphs for each category).
# -----------------------------------------------------------------
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$pTitle = $d.Paragraphs.Last

$pTitle.Range.InsertParagraphAfter()
$pStat = $d.Paragraphs.Last

$pStat.Range.InsertParagraphAfter()
$pBig = $d.Paragraphs.Last

$pBig.Range.InsertParagraphAfter()
$pViz = $d.Paragraphs.Last

$pTitle.Range.Text = "TECHNICAL SKILLS"
$pTitle.Style = "Heading2"

$pStat.Range.Text = "STATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning; Statistical Computing; A/B Testing; Meta-analytical Techniques"

$pBig.Range.Text = "BIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Databases; Data Governance; Streaming Data; Data Pipeline Optimization"

$pViz.Range.Text = "DATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Statistical Reporting; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Business Intelligence; Client Presentation"

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
